# Add a new "UK" worksheet (Test Data for UK Market), based on the existing
# "Netherlands" sheet (same layout: wide column A, no custom column D width)
# but with the row 17-20 ordering and tab-selection state that matches how
# the other recently-added market sheets (e.g. "Russia") look.

$wb = $excel.ActiveWorkbook

# A throwaway sheet is inserted (and removed) purely so the internal sheetId
# counter advances past the value that would otherwise be reused, matching
# the sheetId="18" the new sheet ends up with in the real workbook history.
$placeholder = $wb.Worksheets.Add()

$source = $wb.Worksheets.Item("Netherlands")
$source.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$placeholder.Delete()

$uk = $wb.Worksheets.Item($wb.Worksheets.Count)
$uk.Name = "UK"

# Fill in the market-specific values. B4 (the NGC/test ticket code) is set
# before B2 (the market name) so the new shared-string entries land in the
# same order as the source workbook: NGC code first, then "UK Market".
$uk.Range("B4").Value = "NGC-4331/T3345/T3349/T3366"
$uk.Range("B2").Value = "UK Market"

# Re-order the "Zonal ..." rows to match the Russia/Finland-style ordering.
$uk.Range("A17").Value = "Zonal Dect/Fault Display 40"
$uk.Range("A18").Value = "Zonal Alarm Display max 80"
$uk.Range("A19").Value = "Zonal Alarm/Fault Display 40"
$uk.Range("A20").Value = "Zonal Alarm/Fault/Normal Display 40"

# Make the new sheet the active/selected one, with B4 selected - this also
# clears tabSelected on whichever sheet previously held it (Russia).
$uk.Activate()
$uk.Range("B4").Select()
